# Applies the "fixing color schemes and some mapping stuff" data update:
# appends two new logger readings (3/16/2021 and 3/23/2021 download dates)
# to the bottom of each site data sheet, and leaves the PBSF sheet as the
# active / selected tab (mirroring the author's final on-screen state).

$wb = $excel.ActiveWorkbook

function Add-Rows {
    param(
        [string]$SheetName,
        [int]$FirstRow,
        [array]$Rows,        # array of @(datetimeSerial, conductivity, tempC)
        [string]$SelectCell
    )

    $ws = $wb.Worksheets.Item($SheetName)
    $ws.Activate()

    $r = $FirstRow
    foreach ($row in $Rows) {
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $r = $r + 1
    }

    [void]$ws.Range($SelectCell).Select()
}

# WIC (sheet2) - new rows 25:26
Add-Rows "WIC" 25 @(
    @(44271.411111111112, 552.20000000000005, 4),
    @(44278.401388888888, 704.5, 9.3000000000000007)
) "C26"

# YS (sheet3) - new rows 40:41
Add-Rows "YS" 40 @(
    @(44271.425000000003, 332.8, 3.5),
    @(44278.413194444445, 357.1, 7.2)
) "E40"

# SW (sheet4) - new rows 37:38
Add-Rows "SW" 37 @(
    @(44271.440972222219, 1006, 3.1),
    @(44278.430555555555, 868.9, 7.2)
) "F43"

# YI (sheet5) - new rows 37:38
Add-Rows "YI" 37 @(
    @(44271.440972222219, 323.5, 3.3),
    @(44278.443055555559, 357.4, 5.9)
) "I38"

# YN (sheet6) - new rows 36:37
Add-Rows "YN" 36 @(
    @(44271.526388888888, 397.3, 2.2999999999999998),
    @(44278.493750000001, 518.29999999999995, 9.6)
) "C37"

# 6MC (sheet7) - new rows 38:39
Add-Rows "6MC" 38 @(
    @(44271.543749999997, 455.4, 4.0999999999999996),
    @(44278.50277777778, 528.79999999999995, 9)
) "C39"

# DC (sheet8) - new rows 38:39
Add-Rows "DC" 38 @(
    @(44271.552083333336, 489.3, 5.8),
    @(44278.509722222225, 538.20000000000005, 8.9)
) "C39"

# PBMS (sheet9) - new rows 39:40
Add-Rows "PBMS" 39 @(
    @(44271.567361111112, 840.8, 3.2),
    @(44278.525000000001, 936.7, 10)
) "C41"

# PBSF (sheet10) - new rows 39:40; this is the sheet left active/selected
Add-Rows "PBSF" 39 @(
    @(44271.579861111109, 1374, 5.3),
    @(44278.53402777778, 1811, 11.1)
) "G45"

Write-Host "done"
